$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Borrar puntos al final de frases en vinetas (remove trailing periods
# and stray double comma typo from award descriptions)
$ws.Range("E2").Value = "For ‘trying to quantify the relationship between different countries’ national income inequality and the average amount of mouth-to-mouth kissing’ (Watkins, et al., 2019)"
$ws.Range("E8").Value = "Best overall performance in the MSc"

# Update the saved selection to match the author's cursor position
$ws.Range("E8").Select()
